$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$c_2_D = $ws.Range("D2")
$c_2_D.Value = '29.429.57'
$c_2_E = $ws.Range("E2")
$c_2_E.Value = '  +0.19%  '

# Row 3: 'Ethereum'
$c_3_D = $ws.Range("D3")
$c_3_D.Value = '1.848.18'
$c_3_E = $ws.Range("E3")
$c_3_E.Value = '  +0.28%  '

# Row 4: 'TetherUSD'
$c_4_D = $ws.Range("D4")
$c_4_D.NumberFormat = "@"
$c_4_D.Value = '1.001'
$c_4_D.Style = "Normal"
$c_4_E = $ws.Range("E4")
$c_4_E.Value = '  +0.24%  '

# Row 5: 'BNB'
$c_5_D = $ws.Range("D5")
$c_5_D.NumberFormat = "@"
$c_5_D.Value = '240.73'
$c_5_D.Style = "Normal"
$c_5_E = $ws.Range("E5")
$c_5_E.Value = '  +0.77%  '

# Row 6: 'XRP'
$c_6_D = $ws.Range("D6")
$c_6_D.NumberFormat = "@"
$c_6_D.Value = '0.6275'
$c_6_D.Style = "Normal"
$c_6_E = $ws.Range("E6")
$c_6_E.Value = '  -0.50%  '

# Row 7: 'USDC'
$c_7_D = $ws.Range("D7")
$c_7_D.NumberFormat = "@"
$c_7_D.Value = '1.001'
$c_7_D.Style = "Normal"
$c_7_E = $ws.Range("E7")
$c_7_E.Value = '  +0.13%  '

# Row 8: 'Dogecoin'
$c_8_D = $ws.Range("D8")
$c_8_D.NumberFormat = "@"
$c_8_D.Value = '0.07662'
$c_8_D.Style = "Normal"
$c_8_E = $ws.Range("E8")
$c_8_E.Value = '  +1.71%  '

# Row 9: 'Cardano'
$c_9_D = $ws.Range("D9")
$c_9_D.NumberFormat = "@"
$c_9_D.Value = '0.2914'
$c_9_D.Style = "Normal"
$c_9_E = $ws.Range("E9")
$c_9_E.Value = '  -0.34%  '

# Row 10: 'Solana'
$c_10_D = $ws.Range("D10")
$c_10_D.NumberFormat = "@"
$c_10_D.Value = '24.84'
$c_10_D.Style = "Normal"
$c_10_E = $ws.Range("E10")
$c_10_E.Value = '  +1.92%  '

# Row 11: 'TRON'
$c_11_D = $ws.Range("D11")
$c_11_D.NumberFormat = "@"
$c_11_D.Value = '0.07743'
$c_11_D.Style = "Normal"
$c_11_E = $ws.Range("E11")
$c_11_E.Value = '  +0.47%  '

# Row 12: 'WrappedEther'
$c_12_D = $ws.Range("D12")
$c_12_D.Value = '1.863.83'
$c_12_E = $ws.Range("E12")
$c_12_E.Value = '  -0.78%  '

# Row 13: 'Polkadot'
$c_13_D = $ws.Range("D13")
$c_13_D.NumberFormat = "@"
$c_13_D.Value = '5.027'
$c_13_D.Style = "Normal"
$c_13_E = $ws.Range("E13")
$c_13_E.Value = '  +0.70%  '

# Row 14: 'Polygon'
$c_14_D = $ws.Range("D14")
$c_14_D.NumberFormat = "@"
$c_14_D.Value = '0.6805'
$c_14_D.Style = "Normal"
$c_14_E = $ws.Range("E14")
$c_14_E.Value = '  +0.33%  '

# Row 15: 'ShibaInu'
$c_15_D = $ws.Range("D15")
$c_15_D.NumberFormat = "@"
$c_15_D.Value = '0.00001071'
$c_15_D.Style = "Normal"
$c_15_E = $ws.Range("E15")
$c_15_E.Value = '  +3.76%  '

# Row 16: 'Litecoin'
$c_16_D = $ws.Range("D16")
$c_16_D.NumberFormat = "@"
$c_16_D.Value = '83.42'
$c_16_D.Style = "Normal"

# Row 17: 'Uniswap'
$c_17_D = $ws.Range("D17")
$c_17_D.NumberFormat = "@"
$c_17_D.Value = '6.169'
$c_17_D.Style = "Normal"
$c_17_E = $ws.Range("E17")
$c_17_E.Value = '  +0.25%  '

# Row 18: 'WrappedBTC'
$c_18_D = $ws.Range("D18")
$c_18_D.Value = '29.413.57'
$c_18_E = $ws.Range("E18")
$c_18_E.Value = '  -0.01%  '

# Row 19: 'BitcoinCash'
$c_19_D = $ws.Range("D19")
$c_19_D.NumberFormat = "@"
$c_19_D.Value = '228.34'
$c_19_D.Style = "Normal"
$c_19_E = $ws.Range("E19")
$c_19_E.Value = '  +0.30%  '

# Row 20: 'Avalanche'
$c_20_E = $ws.Range("E20")
$c_20_E.Value = '  -0.32%  '

# Row 21: 'Dai'
$c_21_D = $ws.Range("D21")
$c_21_D.NumberFormat = "@"
$c_21_D.Value = '1.001'
$c_21_D.Style = "Normal"
$c_21_E = $ws.Range("E21")
$c_21_E.Value = '  +0.02%  '

# Row 22: 'Chainlink'
$c_22_D = $ws.Range("D22")
$c_22_D.NumberFormat = "@"
$c_22_D.Value = '7.414'
$c_22_D.Style = "Normal"
$c_22_E = $ws.Range("E22")
$c_22_E.Value = '  -0.45%  '

# Row 23: 'BinanceUSD'
$c_23_D = $ws.Range("D23")
$c_23_D.NumberFormat = "@"
$c_23_D.Value = '1.001'
$c_23_D.Style = "Normal"
$c_23_E = $ws.Range("E23")
$c_23_E.Value = '  +0.02%  '

# Row 24: 'Monero'
$c_24_D = $ws.Range("D24")
$c_24_D.NumberFormat = "@"
$c_24_D.Value = '157.88'
$c_24_D.Style = "Normal"
$c_24_E = $ws.Range("E24")
$c_24_E.Value = '  +0.68%  '

# Row 25: 'Stellar'
$c_25_D = $ws.Range("D25")
$c_25_D.NumberFormat = "@"
$c_25_D.Value = '0.1372'
$c_25_D.Style = "Normal"
$c_25_E = $ws.Range("E25")
$c_25_E.Value = '  -1.32%  '

# Row 26: 'Cosmos'
$c_26_D = $ws.Range("D26")
$c_26_D.NumberFormat = "@"
$c_26_D.Value = '8.397'
$c_26_D.Style = "Normal"
$c_26_E = $ws.Range("E26")
$c_26_E.Value = '  +0.48%  '

# Row 27: 'EthereumClassic'
$c_27_D = $ws.Range("D27")
$c_27_D.NumberFormat = "@"
$c_27_D.Value = '17.67'
$c_27_D.Style = "Normal"
$c_27_E = $ws.Range("E27")
$c_27_E.Value = '  +0.49%  '

# Row 28: 'Toncoin'
$c_28_D = $ws.Range("D28")
$c_28_D.NumberFormat = "@"
$c_28_D.Value = '1.353'
$c_28_D.Style = "Normal"
$c_28_E = $ws.Range("E28")
$c_28_E.Value = '  +6.20%  '

# Row 29: 'PancakeSwap'
$c_29_D = $ws.Range("D29")
$c_29_D.NumberFormat = "@"
$c_29_D.Value = '1.464'
$c_29_D.Style = "Normal"
$c_29_E = $ws.Range("E29")
$c_29_E.Value = '  +0.52%  '

# Row 30: 'Hedera'
$c_30_D = $ws.Range("D30")
$c_30_D.NumberFormat = "@"
$c_30_D.Value = '0.05664'
$c_30_D.Style = "Normal"
$c_30_E = $ws.Range("E30")
$c_30_E.Value = '  +0.78%  '

# Row 31: 'Filecoin'
$c_31_D = $ws.Range("D31")
$c_31_D.NumberFormat = "@"
$c_31_D.Value = '4.117'
$c_31_D.Style = "Normal"
$c_31_E = $ws.Range("E31")
$c_31_E.Value = '  +0.33%  '

# Row 32: 'InternetComputer(DFINITY)'
$c_32_D = $ws.Range("D32")
$c_32_D.NumberFormat = "@"
$c_32_D.Value = '4.033'
$c_32_D.Style = "Normal"
$c_32_E = $ws.Range("E32")
$c_32_E.Value = '  +0.35%  '

# Row 33: 'LidoDAOToken'
$c_33_D = $ws.Range("D33")
$c_33_D.NumberFormat = "@"
$c_33_D.Value = '1.839'
$c_33_D.Style = "Normal"
$c_33_E = $ws.Range("E33")
$c_33_E.Value = '  +0.53%  '

# Row 34: 'ARBITRUM'
$c_34_D = $ws.Range("D34")
$c_34_D.NumberFormat = "@"
$c_34_D.Value = '1.161'
$c_34_D.Style = "Normal"
$c_34_E = $ws.Range("E34")
$c_34_E.Value = '  +0.41%  '

# Row 35: 'ImmutableX'
$c_35_D = $ws.Range("D35")
$c_35_D.NumberFormat = "@"
$c_35_D.Value = '0.7003'
$c_35_D.Style = "Normal"
$c_35_E = $ws.Range("E35")
$c_35_E.Value = '  -0.99%  '

# Row 36: 'HuobiToken'
$c_36_D = $ws.Range("D36")
$c_36_D.NumberFormat = "@"
$c_36_D.Value = '2.588'
$c_36_D.Style = "Normal"
$c_36_E = $ws.Range("E36")
$c_36_E.Value = '  -0.02%  '

# Row 37: 'MXToken'
$c_37_D = $ws.Range("D37")
$c_37_D.NumberFormat = "@"
$c_37_D.Value = '2.762'
$c_37_D.Style = "Normal"
$c_37_E = $ws.Range("E37")
$c_37_E.Value = '  +0.01%  '

# Row 38: 'Maker'
$c_38_D = $ws.Range("D38")
$c_38_D.Value = '1.223.89'
$c_38_E = $ws.Range("E38")
$c_38_E.Value = '  -1.33%  '

# Row 40: 'FraxShare'
$c_40_D = $ws.Range("D40")
$c_40_D.NumberFormat = "@"
$c_40_D.Value = '6.544'
$c_40_D.Style = "Normal"
$c_40_E = $ws.Range("E40")
$c_40_E.Value = '  +3.74%  '

# Row 41: 'TrustWalletToken'
$c_41_D = $ws.Range("D41")
$c_41_D.NumberFormat = "@"
$c_41_D.Value = '0.9022'
$c_41_D.Style = "Normal"
$c_41_E = $ws.Range("E41")
$c_41_E.Value = '  +0.20%  '

# Row 42: 'PaxDollar'
$c_42_E = $ws.Range("E42")
$c_42_E.Value = '  +0.09%  '

# Row 43: 'Quant'
$c_43_B = $ws.Range("B43")
$c_43_B.Value = 'Quant'
$c_43_C = $ws.Range("C43")
$c_43_C.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c_43_D = $ws.Range("D43")
$c_43_D.NumberFormat = "@"
$c_43_D.Value = '101.78'
$c_43_D.Style = "Normal"
$c_43_E = $ws.Range("E43")
$c_43_E.Value = '  -0.19%  '

# Row 44: 'Aave'
$c_44_B = $ws.Range("B44")
$c_44_B.Value = 'Aave'
$c_44_C = $ws.Range("C44")
$c_44_C.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c_44_D = $ws.Range("D44")
$c_44_D.NumberFormat = "@"
$c_44_D.Value = '65.93'
$c_44_D.Style = "Normal"
$c_44_E = $ws.Range("E44")
$c_44_E.Value = '  +0.51%  '

# Row 45: 'BabyDogeCoin'
$c_45_B = $ws.Range("B45")
$c_45_B.Value = 'BabyDogeCoin'
$c_45_C = $ws.Range("C45")
$c_45_C.Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$c_45_D = $ws.Range("D45")
$c_45_D.NumberFormat = "@"
$c_45_D.Value = '0.00000000122'
$c_45_D.Style = "Normal"
$c_45_E = $ws.Range("E45")
$c_45_E.Value = '  +2.32%  '

# Row 46: 'Aptos'
$c_46_D = $ws.Range("D46")
$c_46_D.NumberFormat = "@"
$c_46_D.Value = '7.150'
$c_46_D.Style = "Normal"
$c_46_E = $ws.Range("E46")
$c_46_E.Value = '  +1.23%  '

# Row 47: 'TheSandbox'
$c_47_B = $ws.Range("B47")
$c_47_B.Value = 'TheSandbox'
$c_47_C = $ws.Range("C47")
$c_47_C.Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$c_47_D = $ws.Range("D47")
$c_47_D.NumberFormat = "@"
$c_47_D.Value = '0.4015'
$c_47_D.Style = "Normal"
$c_47_E = $ws.Range("E47")
$c_47_E.Value = '  +0.39%  '

# Row 48: 'EnergySwap'
$c_48_B = $ws.Range("B48")
$c_48_B.Value = 'EnergySwap'
$c_48_C = $ws.Range("C48")
$c_48_C.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c_48_D = $ws.Range("D48")
$c_48_D.NumberFormat = "@"
$c_48_D.Value = '9.020'
$c_48_D.Style = "Normal"
$c_48_E = $ws.Range("E48")
$c_48_E.Value = '  +1.42%  '

# Row 49: 'Algorand'
$c_49_B = $ws.Range("B49")
$c_49_B.Value = 'Algorand'
$c_49_C = $ws.Range("C49")
$c_49_C.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c_49_D = $ws.Range("D49")
$c_49_D.NumberFormat = "@"
$c_49_D.Value = '0.1154'
$c_49_D.Style = "Normal"
$c_49_E = $ws.Range("E49")
$c_49_E.Value = '  +3.23%  '

# Row 50: 'RenderToken'
$c_50_B = $ws.Range("B50")
$c_50_B.Value = 'RenderToken'
$c_50_C = $ws.Range("C50")
$c_50_C.Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c_50_D = $ws.Range("D50")
$c_50_D.NumberFormat = "@"
$c_50_D.Value = '1.673'
$c_50_D.Style = "Normal"
$c_50_E = $ws.Range("E50")
$c_50_E.Value = '  +0.43%  '

# Row 51: 'Cronos'
$c_51_B = $ws.Range("B51")
$c_51_B.Value = 'Cronos'
$c_51_C = $ws.Range("C51")
$c_51_C.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$c_51_D = $ws.Range("D51")
$c_51_D.NumberFormat = "@"
$c_51_D.Value = '0.05718'
$c_51_D.Style = "Normal"
$c_51_E = $ws.Range("E51")
$c_51_E.Value = '  +0.19%  '
